# Accomplishment Log: add "Problem Solved" and "Opportunity Created" columns
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$tbl = $ws.ListObjects.Item(1)

# Add the two new table columns (this grows the table, updates the
# table/worksheet dimension, and appends the header text as new shared
# strings automatically).
$colProblem = $tbl.ListColumns.Add()
$colProblem.Range.Cells.Item(1, 1).Value2 = "Problem Solved"

$colOpportunity = $tbl.ListColumns.Add()
$colOpportunity.Range.Cells.Item(1, 1).Value2 = "Opportunity Created"

# New data-row cells (J2:K2) pick up the same formatting as the other
# "Benefit"/"Person Impacted"/"Business Impact" data cells: vertical
# centering plus word wrap (no horizontal centering).
$newDataCells = $ws.Range("J2:K2")
$newDataCells.VerticalAlignment = -4108
$newDataCells.WrapText = $true

# Every header cell (A1:K1) now shares the same centered + wrapped style
# that used to be reserved for the "Mentioned In Review Period" header.
$headerRow = $ws.Range("A1:K1")
$headerRow.HorizontalAlignment = -4108
$headerRow.WrapText = $true

# Give the newly-visible "Project" column a bit more breathing room, as in
# the authored change.
$ws.Columns.Item(3).ColumnWidth = 15.498697916666666

# Move the selection/view over to show the new columns.
$ws.Range("B1:K1").Select()
$excel.ActiveWindow.ScrollColumn = 2
